$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 378, shifting existing rows 378:475 down to 379:476.
$ws.Rows.Item(378).Insert()

# Populate the newly inserted row 378 with the new data record.
$ws.Range("A378").Value = 5
$ws.Range("B378").Value = "Macroferia Regional de Talca"
$ws.Range("C378").Value = "Maule"
$ws.Range("D378").Value = 45135
$ws.Range("D378").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E378").Value = 7
$ws.Range("F378").Value = 100112009
$ws.Range("G378").Value = "Acelga"
$ws.Range("H378").Value = "Sin especificar"
$ws.Range("I378").Value = "Primera"
$ws.Range("J378").Value = 300
$ws.Range("K378").Value = 1600
$ws.Range("L378").Value = 1600
$ws.Range("M378").Value = 1600
$ws.Range("N378").Value = "`$/docena de atados (4 kilos)"
$ws.Range("O378").Value = "Región del Maule"
$ws.Range("P378").Value = 400
$ws.Range("Q378").Value = 4
$ws.Range("R378").Value = "Hortaliza"
